# Fix the typo "giagnosed" -> "diagnosed" in the sentence on the
# "Dataset" slide (slide 7) and merge the three runs that make up that
# sentence into a single run - matching what PowerPoint produces when a
# user selects the whole sentence and retypes/corrects it as one piece
# of text (the corrected run keeps the formatting of the original first
# run, drops the spell-check "err" flag, and the runs get joined).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

$shp = $null
foreach ($candidate in $s.Shapes) {
    if ($candidate.Name -eq "Content Placeholder 2") {
        $shp = $candidate
    }
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(2)
}

$tr = $shp.TextFrame.TextRange

$fixed = "The sample size was 273 patients of whom 125 were diagnosed with the disease."

$full = $tr.Text
$startMarker = "The sample size"
$endMarker = "with the disease."

$startIdx = $full.IndexOf($startMarker)
$endIdx = $full.IndexOf($endMarker) + $endMarker.Length
$len = $endIdx - $startIdx

$sentence = $tr.Characters($startIdx + 1, $len)
$sentence.Text = $fixed
